# Update "想去人数" (interest count) values on the 展览 and 全部类型 sheets
# to reflect the newly scraped numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F13").Value = 2488
$ws1.Range("F20").Value = 587
$ws1.Range("F25").Value = 2083
$ws1.Range("F32").Value = 2128
$ws1.Range("F39").Value = 723

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 2488
$ws4.Range("F21").Value = 587
$ws4.Range("F26").Value = 2083
$ws4.Range("F33").Value = 2128
$ws4.Range("F40").Value = 723
